$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.05871566666666667
$ws.Range("N2").Value = 0.176147
$ws.Range("O2").Value = 0.008355169877275808
$ws.Range("P2").Value = 0.008355169877275808
$ws.Range("Q2").Value = 0.01269004088966667
$ws.Range("R2").Value = 0.114210368007
$ws.Range("S2").Value = 0.004290111019844772
$ws.Range("T2").Value = 0.004290111019844772

# Row 3
$ws.Range("O3").Value = 0.1868088427899751
$ws.Range("P3").Value = 0.1868088427899751
$ws.Range("S3").Value = 0.0959203327795206
$ws.Range("T3").Value = 0.0959203327795206

# Row 4
$ws.Range("O4").Value = 0.8048359873327491
$ws.Range("P4").Value = 0.8048359873327491
$ws.Range("S4").Value = 0.4132573950189588
$ws.Range("T4").Value = 0.4132573950189588

# Row 5
$ws.Range("M5").Value = 0.05871566666666667
$ws.Range("N5").Value = 0.176147
$ws.Range("O5").Value = 0.008355169877275808
$ws.Range("P5").Value = 0.008355169877275808
$ws.Range("Q5").Value = 0.01202434223288889
$ws.Range("R5").Value = 0.108219080096
$ws.Range("S5").Value = 0.004065058857431037
$ws.Range("T5").Value = 0.004065058857431037

# Row 6
$ws.Range("O6").Value = 0.1868088427899751
$ws.Range("P6").Value = 0.1868088427899751
$ws.Range("S6").Value = 0.09088851001045453
$ws.Range("T6").Value = 0.09088851001045453

# Row 7
$ws.Range("O7").Value = 0.8048359873327491
$ws.Range("P7").Value = 0.8048359873327491
$ws.Range("S7").Value = 0.3915785923137903
$ws.Range("T7").Value = 0.3915785923137903
